$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell-value updates per the cryptos price/volume refresh.
$ws.Range("D2").Value = "58.415.18"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.618.43"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "533.52"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.18"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.567"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.98%  "
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "3.084.32"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "58.343.77"
$ws.Range("E14").Value = "  -1.71%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "20.61"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "2.612.04"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("E18").Value = "  +0.65%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "333.79"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("E20").Value = "  +0.15%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.22"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "66.28"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  -1.44%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("E31").Value = "  +0.53%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "18.74"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "150.23"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.26%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.856"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("E36").Value = "  -1.59%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.812"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("E39").Value = "  +0.82%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "280.01"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -1.23%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "10.67"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0527"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.89"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0934"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "1.934.65"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "17.87"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.87%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.42"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.54%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "112.91"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.69%  "
